# Picture 18 (rId4) on slide 1: apply crop, reposition/resize, and drop the outer shadow.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)

# Crop the image (srcRect l="643" t="936" -> 0.643% / 0.936% of the native image size).
$sh.PictureFormat.CropLeft = 2.1267225
$sh.PictureFormat.CropTop = 2.12706

# New position/size (EMU / 12700 = points).
$sh.Left = 693.0
$sh.Top = 235.82814960629923
$sh.Width = 438.16539370078743
$sh.Height = 300.17192913385827

# Remove the outer shadow effect entirely.
$sh.Shadow.Visible = 0
